$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E8 becomes a formula 12+4+4 (was a plain value 12)
$ws.Range("E8").Formula = "=12+4+4"

# E9 new value
$ws.Range("E9").Value = 3

# E10 new value (new row 10)
$ws.Range("E10").Value = 8

# Update selection to E10 as the active cell
$ws.Range("E10").Select()
